$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-row price / volume updates ---
$ws.Range("D2").Value = "67.175.07"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.476.64"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'585.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'171.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.476.24"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "2.936.70"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "67.032.69"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "2.474.49"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "'351.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").Value = "'4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("E23").Value = "  +0.01%  "

# --- Row 24/25 swap: Litecoin <-> NEARProtocol ---
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").Value = "'4.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'68.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.04%  "

# --- Per-row price / volume updates ---
$ws.Range("D26").Value = "'1.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("D27").Value = "'9.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.63%  "
$ws.Range("D28").Value = "2.571.90"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "0.0₃0906"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "'512.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'162.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "'18.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'1.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").Value = "'143.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "0.0₆0257"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  -0.69%  "
